# Change the table style applied to the table on slide 5 (the "B1- TYPES OF
# FINANCIAL DOCUMENTS" slide) from the custom "Table_0" style to the built-in
# "Themed Style 1 - Accent 1" table style, matching a Table Design ribbon
# gallery pick in PowerPoint.

$p = $ppt.ActivePresentation

$targetSlideIndex = 0
$targetShapeIndex = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $targetSlideIndex = $si
            $targetShapeIndex = $shi
        }
    }
}

$slide = $p.Slides.Item($targetSlideIndex)
$shape = $slide.Shapes.Item($targetShapeIndex)
$table = $shape.Table

$table.ApplyStyle("{D449F14A-2BBA-474C-932B-564049EC565D}")
